$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") values for rows 2-267 all change from 45189 to 45190
# (the date serial number advances by one day, 2023-09-20 -> 2023-09-21).
$rng = $ws.Range("C2:C267")
$rng.Value2 = 45190
